# Apply the edit described by the commit:
#   "Code added for Zip file compressor and is working fine"
#
# Underlying data change: the "SK-Moving Leave" test row (which was Failing)
# is removed, and the remaining rows' TestResult flips from Failed -> Passed
# (i.e. the zip-compressor test run now reports success).
#
#  - Delete row 3 ("SK-Moving Leave" / 25/04/2025) entirely; row 4
#    ("SK-Doctor Visit Family Member Care Full") shifts up to become row 3.
#  - Change the "TestResult" column value from "Failed" to "Passed" for the
#    remaining data rows.
#  - Restore the print/zoom settings and selection left in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 3 (SK-Moving Leave), shifting row 4 (SK-Doctor Visit...) up into row 3.
$ws.Rows.Item(3).Delete()

# Update TestResult column (G) values from "Failed" to "Passed" for the remaining data rows.
$ws.Range("G2:G3").Value = "Passed"

# Match the saved view/print state.
$excel.ActiveWindow.Zoom = 100
$ws.PageSetup.Zoom = 100

# Selection moved to B12 in the saved file.
$ws.Range("B12").Select()
